$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.434592
$ws.Range("H2").Value = 58.303776
$ws.Range("I2").Value = 0.1244167820899015
$ws.Range("J2").Value = 0.1244167820899015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 459.3520842462186
$ws.Range("R2").Value = 4134.168758215968
$ws.Range("S2").Value = 0.008494984522017773
$ws.Range("T2").Value = 0.008494984522017772
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.434592
$ws.Range("H3").Value = 58.303776
$ws.Range("I3").Value = 0.1244167820899015
$ws.Range("J3").Value = 0.1244167820899015
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 3523.265919899669
$ws.Range("R3").Value = 31709.39327909702
$ws.Range("S3").Value = 0.06515718657424768
$ws.Range("T3").Value = 0.06515718657424766
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.434592
$ws.Range("H4").Value = 58.303776
$ws.Range("I4").Value = 0.1244167820899015
$ws.Range("J4").Value = 0.1244167820899015
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 2159.193788965429
$ws.Range("R4").Value = 19432.74410068886
$ws.Range("S4").Value = 0.03993084704817952
$ws.Range("T4").Value = 0.03993084704817952
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.434592
$ws.Range("H5").Value = 58.303776
$ws.Range("I5").Value = 0.1244167820899015
$ws.Range("J5").Value = 0.1244167820899015
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 585.8176710832853
$ws.Range("R5").Value = 5272.359039749567
$ws.Range("S5").Value = 0.01083376394545656
$ws.Range("T5").Value = 0.01083376394545656
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 48.891945
$ws.Range("H6").Value = 146.675835
$ws.Range("I6").Value = 0.3129974875220664
$ws.Range("J6").Value = 0.3129974875220664
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 1155.600119549795
$ws.Range("R6").Value = 10400.40107594816
$ws.Range("S6").Value = 0.02137098201116567
$ws.Range("T6").Value = 0.02137098201116567
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 48.891945
$ws.Range("H7").Value = 146.675835
$ws.Range("I7").Value = 0.3129974875220664
$ws.Range("J7").Value = 0.3129974875220664
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 8863.542058207811
$ws.Range("R7").Value = 79771.8785238703
$ws.Range("S7").Value = 0.1639170805511562
$ws.Range("T7").Value = 0.1639170805511562
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.891945
$ws.Range("H8").Value = 146.675835
$ws.Range("I8").Value = 0.3129974875220664
$ws.Range("J8").Value = 0.3129974875220664
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 5431.921800799285
$ws.Range("R8").Value = 48887.29620719356
$ws.Range("S8").Value = 0.1004547344077512
$ws.Range("T8").Value = 0.1004547344077512
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.891945
$ws.Range("H9").Value = 146.675835
$ws.Range("I9").Value = 0.3129974875220664
$ws.Range("J9").Value = 0.3129974875220664
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 1473.75182121817
$ws.Range("R9").Value = 13263.76639096353
$ws.Range("S9").Value = 0.02725469055199333
$ws.Range("T9").Value = 0.02725469055199332
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 65.19353
$ws.Range("H10").Value = 195.58059
$ws.Range("I10").Value = 0.4173573191390618
$ws.Range("J10").Value = 0.4173573191390618
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 1540.901084255763
$ws.Range("R10").Value = 13868.10975830187
$ws.Range("S10").Value = 0.02849650912587727
$ws.Range("T10").Value = 0.02849650912587727
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 65.19353
$ws.Range("H11").Value = 195.58059
$ws.Range("I11").Value = 0.4173573191390618
$ws.Range("J11").Value = 0.4173573191390618
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 11818.83017904141
$ws.Range("R11").Value = 106369.4716113727
$ws.Range("S11").Value = 0.2185704231734741
$ws.Range("T11").Value = 0.218570423173474
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 65.19353
$ws.Range("H12").Value = 195.58059
$ws.Range("I12").Value = 0.4173573191390618
$ws.Range("J12").Value = 0.4173573191390618
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 7243.036800398556
$ws.Range("R12").Value = 65187.331203587
$ws.Range("S12").Value = 0.1339484191363988
$ws.Range("T12").Value = 0.1339484191363988
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 65.19353
$ws.Range("H13").Value = 195.58059
$ws.Range("I13").Value = 0.4173573191390618
$ws.Range("J13").Value = 0.4173573191390618
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 1965.131139068846
$ws.Range("R13").Value = 17686.18025161962
$ws.Range("S13").Value = 0.03634196770331172
$ws.Range("T13").Value = 0.03634196770331172
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 22.685484
$ws.Range("H14").Value = 68.05645200000001
$ws.Range("I14").Value = 0.1452284112489703
$ws.Range("J14").Value = 0.1452284112489703
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 536.189509794404
$ws.Range("R14").Value = 4825.705588149636
$ws.Range("S14").Value = 0.00991597021715104
$ws.Range("T14").Value = 0.00991597021715104
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 22.685484
$ws.Range("H15").Value = 68.05645200000001
$ws.Range("I15").Value = 0.1452284112489703
$ws.Range("J15").Value = 0.1452284112489703
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 4112.614900978073
$ws.Range("R15").Value = 37013.53410880265
$ws.Range("S15").Value = 0.07605625646862618
$ws.Range("T15").Value = 0.07605625646862617
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 22.685484
$ws.Range("H16").Value = 68.05645200000001
$ws.Range("I16").Value = 0.1452284112489703
$ws.Range("J16").Value = 0.1452284112489703
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 2520.369666236092
$ws.Range("R16").Value = 22683.32699612483
$ws.Range("S16").Value = 0.0466102191297828
$ws.Range("T16").Value = 0.0466102191297828
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 22.685484
$ws.Range("H17").Value = 68.05645200000001
$ws.Range("I17").Value = 0.1452284112489703
$ws.Range("J17").Value = 0.1452284112489703
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 683.809436507704
$ws.Range("R17").Value = 6154.284928569336
$ws.Range("S17").Value = 0.01264596543341026
$ws.Range("T17").Value = 0.01264596543341026
